# Update the dSF (column F) values for the rows identified in the diff.
# These edits reflect a "repull data, push all data, mean calculation" update
# where column F (dSF) values were recalculated/changed for a subset of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    5  = -4
    19 = -2
    25 = 2
    26 = 1
    27 = -1
    38 = -2
    40 = 3
    44 = 1
    52 = 1
    54 = -1
    64 = 1
    68 = 0
    70 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
